$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted at row 200 (pushing the
# existing rows 200-253 down to 201-254). Insert a fresh row above the
# current row 200 so everything below cascades down by one, matching the
# "insert at top of block" pattern seen in the rest of the sheet.
$ws.Rows(200).Insert()

# Populate the newly inserted row 200 with the new observation. Columns
# that are identical to the row that used to occupy this position
# (A,B,C,E,F,G,H,I,N,O,Q,R) are simply restated; D,J,K,L,M,P carry the
# new values.
$ws.Cells.Item(200, 1).Value2 = 11
$ws.Cells.Item(200, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(200, 3).Value2 = "Bíobío"
$ws.Cells.Item(200, 4).Value2 = 45135
$ws.Cells.Item(200, 5).Value2 = 8
$ws.Cells.Item(200, 6).Value2 = 100112032
$ws.Cells.Item(200, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(200, 8).Value2 = "Sin especificar"
$ws.Cells.Item(200, 9).Value2 = "Primera"
$ws.Cells.Item(200, 10).Value2 = 100
$ws.Cells.Item(200, 11).Value2 = 15000
$ws.Cells.Item(200, 12).Value2 = 16000
$ws.Cells.Item(200, 13).Value2 = 15500
$ws.Cells.Item(200, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(200, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(200, 16).Value2 = 310
$ws.Cells.Item(200, 17).Value2 = 50
$ws.Cells.Item(200, 18).Value2 = "Hortaliza"
